$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("110_1")

# Swap the "Total new nominations" (B34) and "Total carryover nominations" (B35)
# values, and give B34 the same number format (thousands separator) as B35.
$ws.Range("B34").Value = 23147
$ws.Range("B34").NumberFormat = $ws.Range("B35").NumberFormat

$ws.Range("B35").Value = 0
